$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for the table body (rows 2-9), columns A-E.
$data = @(
  @("Bentleigh",  "Coles - Bentleigh  5/7 Vickery St, Bentleigh", "24/12/20 12:30pm-1:00pm", "Case shopped", "old"),
  @("Chadstone",  "Chadstone Shopping Centre, Princes Highway", "26/12/20 6am-1:30pm", "Case did not attend during infectious period but may have acquired their illness here. Case visited Culture Kings, Huffer, JD Sports, Jay Jays, H&M, Uniqlo, Myer, Superdry, Footlocker, Dumplings Plus", "old"),
  @("Chadstone",  "Chadstone Shopping Centre, Princes Highway", "26/12/20 6am-1:30pm", "Case did not attend during infectious period but may have acquired their illness here. If you attended Chadstone Shopping Centre but did not attend to any of the acquisition site stores listed above, monitor for symptoms - If symptoms develop, immediately get tested and isolate until you receive a negative result.", "new"),
  @("Hallam",     "Coles Hallam, 2 Princes Domain Drive, Hallam, VIC 3803", "30/12/20 6:15am - 6:30am", "Case shopped in store", "new"),
  @("Moorabbin",  "Costco Moorabbin  8 Chifley Drive, Moorabbin Airport VIC 3194", "30/12/20 10:45am - 12:15pm and 4:00pm- 5:50pm", "Case shopped in store", "new"),
  @("Moorabbin",  "Costco Moorabbin  8 Chifley Drive, Moorabbin Airport VIC 3194", "30/12/20 4:00pm- 5:50pm", "Case shopped in store", "old"),
  @("Mordialloc", "Woodlands Golf Club, 109 White Street, Mordialloc, VIC 3195", "28/12/20 12:00pm - 6:00pm", "Case attended course", "new"),
  @("Wonthaggi",  "Wonthaggi Plaza Shopping Centre, 2 Biggs Drive, Wonthaggi, VIC 3995", "28/12/20 1:30pm - 2:30pm", "Kmart - shopped for 15 mins", "new")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $startRow + $i
  $rowVals = $data[$i]
  for ($c = 0; $c -lt $rowVals.Length; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
  }
}

$ws.Columns.Item(1).ColumnWidth = 9.46484375
$ws.Columns.Item(2).ColumnWidth = 56.73046875
$ws.Columns.Item(3).ColumnWidth = 41.86328125
$ws.Columns.Item(4).ColumnWidth = 27.73046875
$ws.Columns.Item(5).ColumnWidth = 4.46484375
